# Updating with new SIRFER data
# Appends a new block of "utah_untreated_30" treatment measurements
# (samples BOO-11, FE01-FE08, TIE) to the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of raw interlab data (sample, d13C, d18O) for treatment "utah_untreated_30"
$newData = @(
    @("BOO-11", -14.549566,              -8.5316145999999993),
    @("FE01",   -9.3290600000000001,     -1.033345),
    @("FE02",   -15.508796999999999,     -7.5520424999999998),
    @("FE03",   -8.2918199999999995,     -3.8960571000000002),
    @("FE04",   -10.58713,               -1.5984670999999999),
    @("FE05",   -15.971700999999999,     -5.6693895999999997),
    @("FE06",   -13.74333,               -4.2283524999999997),
    @("FE07",   -10.55585,                0.5360066),
    @("FE08",   -9.4439240000000009,     -0.81820850000000001),
    @("TIE",    -15.208564000000001,     -8.2727003000000003)
)

$startRow = 102
$treatment = "utah_untreated_30"

# Write column A (sample) first for every new row so that the brand-new
# shared-string entries are introduced in sample order (FE01..FE08)
# before the repeated treatment label is introduced.
$row = $startRow
foreach ($d in $newData) {
    $ws.Cells.Item($row, 1).Value2 = $d[0]
    $row = $row + 1
}

# Then write column C (treatment) for every new row.
$row = $startRow
foreach ($d in $newData) {
    $ws.Cells.Item($row, 3).Value2 = $treatment
    $row = $row + 1
}

# Finally write the numeric measurement columns B (d13C) and D (d18O).
$row = $startRow
foreach ($d in $newData) {
    $ws.Cells.Item($row, 2).Value2 = $d[1]
    $ws.Cells.Item($row, 4).Value2 = $d[2]
    $row = $row + 1
}

# Update the sheet view: scroll down toward the newly added rows and
# select the newly populated treatment column.
[void]$ws.Range("C103:C111").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1

# Force an explicit page orientation (portrait) to be written out.
$ws.PageSetup.Orientation = 1
